# Create the finance statement report: add the missing student record for
# "tavoz" (Student ID CIC202307009) as a new row in the students table.
#
# The new row is inserted at row 6, which shifts every following row
# (previously rows 6-10) down by one (to rows 7-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 6, pushing the existing rows down.
$ws.Rows.Item(6).Insert()

# Fill in the new student's finance record.
$ws.Cells.Item(6, 1).Value = "tavoz"          # Name
$ws.Cells.Item(6, 2).Value = "mafura"         # Surname
$ws.Cells.Item(6, 3).Value = "CIC202307009"   # Student ID
$ws.Cells.Item(6, 4).Value = 15.0             # Fees
$ws.Cells.Item(6, 5).Value = 15.0             # Owed Fees
$ws.Cells.Item(6, 6).Value = 0.0              # Paid Amount
